$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''51.775.34'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '''2.817.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''355.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.31%  '
$ws.Range('D6').Value = '''111.50'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.45%  '
$ws.Range('D7').Value = '''0.565'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.07%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '''0.596'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.25%  '
$ws.Range('D10').Value = '''40.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.09%  '
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = '''0.132'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').Value = '''19.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('E14').Value = '  +1.27%  '
$ws.Range('D15').Value = '''3.260.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.56%  '
$ws.Range('D16').Value = '''2.807.43'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '''0.917'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.82%  '
$ws.Range('D18').Value = '''51.657.52'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').Value = '''7.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.01%  '
$ws.Range('D20').Value = '''3.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.54%  '
$ws.Range('D21').Value = '''13.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.80%  '
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('D23').Value = '''69.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = '''267.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.19%  '
$ws.Range('D25').Value = '''2.79'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('E26').Value = '  +1.33%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('E30').Value = '  +27.52%  '
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').Value = '''52.63'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.04%  '
$ws.Range('D33').Value = '''34.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  +2.28%  '
$ws.Range('D35').Value = '''5.42'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.43%  '
$ws.Range('E36').Value = '  +2.69%  '
$ws.Range('E38').Value = '  +1.78%  '
$ws.Range('E39').Value = '  -3.36%  '
$ws.Range('E40').Value = '  -3.59%  '
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('D42').Value = '''2.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.48%  '
$ws.Range('D43').Value = '''23.19'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('D44').Value = '''125.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('E45').Value = '  -2.19%  '
$ws.Range('D46').Value = '''2.094.56'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('D47').Value = '''3.33'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.61%  '
$ws.Range('E49').Value = '  +7.75%  '
$ws.Range('D50').Value = '''0.987'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.38%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').Value = '''9.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.92%  '
